$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("이채현")

# Row 4: 신청 메소드 작성
$ws.Range("A4").Value = "신청 메소드 작성"
$ws.Range("B4").Value = "신청 기능의 메소드를 작성"
$ws.Range("C4").Value = "2019-05-17"
$ws.Range("D4").Value = "2019-05-17"
$ws.Range("E4").Value = "완료"

# Row 5: 선발 메소드 작성
$ws.Range("A5").Value = "선발 메소드 작성"
$ws.Range("B5").Value = "선발 기능의 메소드 작성"
$ws.Range("C5").Value = "2019-05-17"
$ws.Range("D5").Value = "2019-05-17"
$ws.Range("E5").Value = "반만 완료"
$ws.Range("F5").Value = "detail 구현 필요"

# Row 6: 클래스 구조 변경
$ws.Range("A6").Value = "클래스 구조 변경"
$ws.Range("B6").Value = "information, operation 소유관계 변경"
$ws.Range("C6").Value = "2019-05-17"
$ws.Range("D6").Value = "2019-05-17"
$ws.Range("E6").Value = "완료"

# Row 7: 기타 클래스 정리
$ws.Range("A7").Value = "기타 클래스 정리"
$ws.Range("B7").Value = "0516 회의 결과 정리된 클래스 구조로 변경 (소유관계 변경 부분 제외)"
$ws.Range("C7").Value = "2019-05-17"
$ws.Range("D7").Value = "2019-05-17"
$ws.Range("E7").Value = "완료"

[void]$ws.Range("B7").Select()
